# Auto-generated Excel COM-interop script
# Applies updated FFXIV leve-crafting profit figures across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the upstream market-data refresh.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4081.2222
$ws.Range("I2").Value = 1349.2
$ws.Range("K2").Value = 1349.2
$ws.Range("M2").Value = -1236.2
$ws.Range("H43").Value = 4500
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5138
$ws.Range("H76").Value = 5500
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 5500
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H96").Value = 19564.637
$ws.Range("I96").Value = 23106.777
$ws.Range("K96").Value = 69320.33099999999
$ws.Range("M96").Value = -67947.33099999999
$ws.Range("H112").Value = 1280.5186
$ws.Range("I112").Value = 1250
$ws.Range("J112").Value = 1282.96
$ws.Range("K112").Value = 3750
$ws.Range("L112").Value = 3848.88
$ws.Range("M112").Value = -2642
$ws.Range("N112").Value = -6064.88
$ws.Range("H113").Value = 4150.5
$ws.Range("I113").Value = 3976
$ws.Range("K113").Value = 3976
$ws.Range("M113").Value = -722
$ws.Range("H135").Value = 1925
$ws.Range("I135").Value = 1415
$ws.Range("K135").Value = 12735
$ws.Range("M135").Value = -10200
$ws.Range("H138").Value = 3264.99
$ws.Range("J138").Value = 3266.4675
$ws.Range("L138").Value = 9799.4025
$ws.Range("N138").Value = -20079.4025

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2581.1
$ws.Range("I63").Value = 2211.8572
$ws.Range("J63").Value = 3442.6667
$ws.Range("K63").Value = 2211.8572
$ws.Range("L63").Value = 3442.6667
$ws.Range("M63").Value = -1525.8572
$ws.Range("N63").Value = -4814.6667
$ws.Range("H66").Value = 2581.1
$ws.Range("I66").Value = 2211.8572
$ws.Range("J66").Value = 3442.6667
$ws.Range("K66").Value = 11059.286
$ws.Range("L66").Value = 17213.3335
$ws.Range("M66").Value = -7627.286
$ws.Range("N66").Value = -24077.3335
$ws.Range("H122").Value = 2090
$ws.Range("I122").Value = 1725
$ws.Range("J122").Value = 3550
$ws.Range("K122").Value = 5175
$ws.Range("L122").Value = 10650
$ws.Range("M122").Value = -2725
$ws.Range("N122").Value = -15550
$ws.Range("H127").Value = 40000
$ws.Range("I127").Value = 40000
$ws.Range("K127").Value = 40000
$ws.Range("M127").Value = -35040
$ws.Range("H132").Value = 1810.6451
$ws.Range("J132").Value = 2151.1667
$ws.Range("L132").Value = 6453.500100000001
$ws.Range("N132").Value = -11513.5001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 9999.5
$ws.Range("I23").Value = 9999
$ws.Range("K23").Value = 9999
$ws.Range("M23").Value = -9716
$ws.Range("H86").Value = 2044.5
$ws.Range("I86").Value = 2044.5
$ws.Range("K86").Value = 2044.5
$ws.Range("M86").Value = -921.5
$ws.Range("H89").Value = 2044.5
$ws.Range("I89").Value = 2044.5
$ws.Range("K89").Value = 10222.5
$ws.Range("M89").Value = -4606.5
$ws.Range("H105").Value = 4085.5
$ws.Range("I105").Value = 3954.8572
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3954.8572
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -2207.8572
$ws.Range("N105").Value = -8494

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 11095.869
$ws.Range("J22").Value = 35945.43
$ws.Range("L22").Value = 35945.43
$ws.Range("N22").Value = -36645.43
$ws.Range("H111").Value = 25500
$ws.Range("J111").Value = 25500
$ws.Range("L111").Value = 25500
$ws.Range("N111").Value = -33680
$ws.Range("H134").Value = 3132.45
$ws.Range("I134").Value = 2953.3125
$ws.Range("J134").Value = 3849
$ws.Range("K134").Value = 8859.9375
$ws.Range("L134").Value = 11547
$ws.Range("M134").Value = -6324.9375
$ws.Range("N134").Value = -16617
$ws.Range("H135").Value = 96926.336
$ws.Range("J135").Value = 96926.336
$ws.Range("L135").Value = 96926.336
$ws.Range("N135").Value = -107066.336

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4499.5
$ws.Range("I14").Value = 4499.5
$ws.Range("K14").Value = 13498.5
$ws.Range("M14").Value = -13325.5
$ws.Range("H38").Value = 1915.7858
$ws.Range("I38").Value = 1401.5
$ws.Range("K38").Value = 4204.5
$ws.Range("M38").Value = -3857.5
$ws.Range("H68").Value = 15630700
$ws.Range("J68").Value = 31260500
$ws.Range("L68").Value = 93781500
$ws.Range("N68").Value = -93783122
$ws.Range("H71").Value = 15630700
$ws.Range("J71").Value = 31260500
$ws.Range("L71").Value = 281344500
$ws.Range("N71").Value = -281352612
$ws.Range("H129").Value = 2537.8
$ws.Range("J129").Value = 6500
$ws.Range("L129").Value = 19500
$ws.Range("N129").Value = -29500
$ws.Range("H131").Value = 3038.2307
$ws.Range("J131").Value = 3349.9
$ws.Range("L131").Value = 10049.7
$ws.Range("N131").Value = -20129.7
$ws.Range("H139").Value = 3178.4443
$ws.Range("I139").Value = 2950.75
$ws.Range("K139").Value = 8852.25
$ws.Range("M139").Value = -3712.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5620.125
$ws.Range("J70").Value = 5995.2
$ws.Range("L70").Value = 5995.2
$ws.Range("N70").Value = -6535.2
$ws.Range("H73").Value = 5620.125
$ws.Range("J73").Value = 5995.2
$ws.Range("L73").Value = 5995.2
$ws.Range("N73").Value = -7867.2
$ws.Range("H97").Value = 278.9091
$ws.Range("I97").Value = 291.77777
$ws.Range("K97").Value = 291.77777
$ws.Range("M97").Value = 204.22223
$ws.Range("H122").Value = 2113.4167
$ws.Range("I122").Value = 795.75
$ws.Range("J122").Value = 4748.75
$ws.Range("K122").Value = 2387.25
$ws.Range("L122").Value = 14246.25
$ws.Range("M122").Value = 62.75
$ws.Range("N122").Value = -19146.25
$ws.Range("H132").Value = 1437.5
$ws.Range("I132").Value = 1437.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4312.5
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("M132").Value = -1782.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8743.950000000001
$ws.Range("I7").Value = 7897.75
$ws.Range("J7").Value = 8955.5
$ws.Range("K7").Value = 7897.75
$ws.Range("L7").Value = 8955.5
$ws.Range("M7").Value = -7785.75
$ws.Range("N7").Value = -9179.5
$ws.Range("H13").Value = 1347
$ws.Range("I13").Value = 1183.75
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 1183.75
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = -1043.75
$ws.Range("N13").Value = -2280
$ws.Range("H22").Value = 2728.25
$ws.Range("I22").Value = 3144
$ws.Range("K22").Value = 3144
$ws.Range("M22").Value = -2849
$ws.Range("H27").Value = 2728.25
$ws.Range("I27").Value = 3144
$ws.Range("K27").Value = 3144
$ws.Range("M27").Value = -3037
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828
$ws.Range("H46").Value = 26875.098
$ws.Range("J46").Value = 3024.6667
$ws.Range("L46").Value = 3024.6667
$ws.Range("N46").Value = -3400.6667
$ws.Range("H126").Value = 8743.950000000001
$ws.Range("I126").Value = 7897.75
$ws.Range("J126").Value = 8955.5
$ws.Range("K126").Value = 23693.25
$ws.Range("L126").Value = 26866.5
$ws.Range("M126").Value = -21223.25
$ws.Range("N126").Value = -31806.5
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H132").Value = 4555.533
$ws.Range("I132").Value = 3985.0908
$ws.Range("J132").Value = 6124.25
$ws.Range("K132").Value = 11955.2724
$ws.Range("L132").Value = 18372.75
$ws.Range("M132").Value = -9425.2724
$ws.Range("N132").Value = -23432.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1583.3334
$ws.Range("J81").Value = 1750
$ws.Range("L81").Value = 3500
$ws.Range("N81").Value = -5622
$ws.Range("H84").Value = 1583.3334
$ws.Range("J84").Value = 1750
$ws.Range("L84").Value = 17500
$ws.Range("N84").Value = -28108
$ws.Range("H122").Value = 4498.5713
$ws.Range("I122").Value = 3798.1
$ws.Range("K122").Value = 11394.3
$ws.Range("M122").Value = -8944.299999999999
